# Updating filtered feeds from workflow
# Appends three new rows to the "Filtered Feeds" worksheet describing the
# new Guardant / Merck companion-diagnostics story (one row per source
# link: GenomeWeb, 360Dx, FierceBiotech), mirroring the existing layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$hyperlinkStyle = $ws.Range("A2").Style

$keywords = "companion diagnostics"
$title = "Guardant, Merck Partner to Develop and Commercialize Companion Diagnostics"
$titleFierce = '<a href="https://www.fiercebiotech.com/medtech/guardant-health-pens-merck-co-cancer-collab-next-gen-tests-and-companion-diagnostics" hreflang="en">Guardant Health pens Merck &amp; Co. cancer collab for next-gen tests, companion diagnostics </a>'

$rows = @(
    @{ Row = 80; Link = "https://www.genomeweb.com/cancer/guardant-merck-partner-develop-and-commercialize-companion-diagnostics"; Title = $title },
    @{ Row = 81; Link = "https://www.360dx.com/cancer/guardant-merck-partner-develop-and-commercialize-companion-diagnostics"; Title = $title },
    @{ Row = 82; Link = "https://www.fiercebiotech.com/medtech/guardant-health-pens-merck-co-cancer-collab-next-gen-tests-and-companion-diagnostics"; Title = $titleFierce }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $linkCell = $ws.Cells.Item($rowNum, 1)

    $linkCell.Value = $r.Link
    $ws.Hyperlinks.Add($linkCell, $r.Link)
    $linkCell.Style = $hyperlinkStyle

    $ws.Cells.Item($rowNum, 2).Value = $keywords
    $ws.Cells.Item($rowNum, 3).Value = $r.Title
}
